$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so numeric-looking strings
# like "587.78" are not auto-converted to floating point numbers by Excel,
# matching the original inline-string representation used in the workbook.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "66.351.40"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "3.459.89"
$ws.Range("E3").Value = "  -1.92%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "587.78"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").Value = "176.14"
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("D7").Value = "0.612"
$ws.Range("E7").Value = "  +1.60%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "3.460.09"
$ws.Range("E9").Value = "  -1.91%  "
$ws.Range("E10").Value = "  -2.25%  "
$ws.Range("D11").Value = "6.96"
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("D12").Value = "0.417"
$ws.Range("E12").Value = "  -1.86%  "
$ws.Range("D13").Value = "4.061.39"
$ws.Range("E13").Value = "  -1.87%  "
$ws.Range("E14").Value = "  +1.65%  "
$ws.Range("D15").Value = "30.06"
$ws.Range("E15").Value = "  -2.16%  "
$ws.Range("D16").Value = "66.289.97"
$ws.Range("E16").Value = "  -0.89%  "
$ws.Range("E17").Value = "  -1.55%  "
$ws.Range("D18").Value = "3.460.96"
$ws.Range("E18").Value = "  -1.92%  "
$ws.Range("E19").Value = "  -2.35%  "
$ws.Range("E20").Value = "  -2.17%  "
$ws.Range("D21").Value = "373.96"
$ws.Range("E21").Value = "  -2.34%  "
$ws.Range("D22").Value = "7.60"
$ws.Range("E22").Value = "  -3.54%  "
$ws.Range("D23").Value = "73.42"
$ws.Range("E23").Value = "  +2.40%  "
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "0.0000126"
$ws.Range("E25").Value = "  +3.03%  "
$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").Value = "0.536"
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "9.92"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "0.178"
$ws.Range("E28").Value = "  +2.25%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "5.85"
$ws.Range("E30").Value = "  -2.83%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "2.00"
$ws.Range("E31").Value = "  -0.99%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "23.72"
$ws.Range("E32").Value = "  -3.95%  "
$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "1.28"
$ws.Range("E34").Value = "  -7.16%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").Value = "7.02"
$ws.Range("E35").Value = "  -3.35%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "1.55"
$ws.Range("E36").Value = "  -1.45%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "161.33"
$ws.Range("E37").Value = "  +1.33%  "
$ws.Range("B38").Value = "Mantle"
$ws.Range("C38").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D38").Value = "0.884"
$ws.Range("E38").Value = "  -0.93%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").Value = "28.47"
$ws.Range("E39").Value = "  -3.31%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "1.81"
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "2.62"
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "4.50"
$ws.Range("E42").Value = "  -0.88%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.763.74"
$ws.Range("E43").Value = "  +1.12%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "6.44"
$ws.Range("E44").Value = "  -3.01%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "0.0693"
$ws.Range("E45").Value = "  -2.16%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "25.38"
$ws.Range("E46").Value = "  -1.40%  "
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").Value = "338.22"
$ws.Range("E47").Value = "  +2.37%  "
$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D48").Value = "40.02"
$ws.Range("E48").Value = "  -1.62%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "0.0292"
$ws.Range("E49").Value = "  -2.61%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "0.103"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D51").Value = "0.994"
$ws.Range("E51").Value = "  -3.36%  "

# Restore the default cell style so no stray style index is left on the cells
# (keeps formatting identical to the original, unstyled data cells).
$priceRange.Style = "Normal"

